# Auto-generated Excel COM-interop script
# Applies profit-recalculation updates across the Jenova_Profits workbook's job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 127366.375
$ws.Range("J28").Value = 349
$ws.Range("L28").Value = 349
$ws.Range("N28").Value = -1319
$ws.Range("H42").Value = 70.14286
$ws.Range("J42").Value = 93
$ws.Range("L42").Value = 279
$ws.Range("N42").Value = -739
$ws.Range("H98").Value = 3310.7307
$ws.Range("I98").Value = 2925.3333
$ws.Range("J98").Value = 4929.4
$ws.Range("K98").Value = 2925.3333
$ws.Range("L98").Value = 4929.4
$ws.Range("M98").Value = -1427.3333
$ws.Range("N98").Value = -7925.4
$ws.Range("H100").Value = 9698.941000000001
$ws.Range("I100").Value = 1713.1428
$ws.Range("K100").Value = 1713.1428
$ws.Range("M100").Value = -1172.1428
$ws.Range("H122").Value = 3310.7307
$ws.Range("I122").Value = 2925.3333
$ws.Range("J122").Value = 4929.4
$ws.Range("K122").Value = 8775.999899999999
$ws.Range("L122").Value = 14788.2
$ws.Range("M122").Value = -6325.999899999999
$ws.Range("N122").Value = -19688.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35333.332
$ws.Range("I2").Value = 50500
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 50500
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -50387
$ws.Range("N2").Value = -5226
$ws.Range("H32").Value = 3602.3433
$ws.Range("I32").Value = 3318.6094
$ws.Range("J32").Value = 9655.333000000001
$ws.Range("K32").Value = 3318.6094
$ws.Range("L32").Value = 9655.333000000001
$ws.Range("M32").Value = -3031.6094
$ws.Range("N32").Value = -10229.333
$ws.Range("H45").Value = 2481.1538
$ws.Range("I45").Value = 1875.5
$ws.Range("K45").Value = 1875.5
$ws.Range("M45").Value = -1498.5
$ws.Range("H61").Value = 2130.6956
$ws.Range("I61").Value = 1523.4286
$ws.Range("J61").Value = 8507
$ws.Range("K61").Value = 1523.4286
$ws.Range("L61").Value = 8507
$ws.Range("M61").Value = -1311.4286
$ws.Range("N61").Value = -8931
$ws.Range("H74").Value = 1185.25
$ws.Range("I74").Value = 1164.2667
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1164.2667
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -290.2666999999999
$ws.Range("N74").Value = -3248
$ws.Range("H77").Value = 1185.25
$ws.Range("I77").Value = 1164.2667
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 5821.3335
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -1453.3335
$ws.Range("N77").Value = -16236
$ws.Range("H110").Value = 187434.62
$ws.Range("I110").Value = 187434.62
$ws.Range("K110").Value = 187434.62
$ws.Range("M110").Value = -185389.62
$ws.Range("H116").Value = 35333.332
$ws.Range("I116").Value = 50500
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 50500
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -48206
$ws.Range("N116").Value = -9588
$ws.Range("H132").Value = 3139.8572
$ws.Range("I132").Value = 3468.9644
$ws.Range("K132").Value = 10406.8932
$ws.Range("M132").Value = -7876.893199999999
$ws.Range("H136").Value = 2130.6956
$ws.Range("I136").Value = 1523.4286
$ws.Range("J136").Value = 8507
$ws.Range("K136").Value = 4570.2858
$ws.Range("L136").Value = 25521
$ws.Range("M136").Value = -2020.2858
$ws.Range("N136").Value = -30621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35333.332
$ws.Range("I3").Value = 50500
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 50500
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -50386
$ws.Range("N3").Value = -5228
$ws.Range("H80").Value = 1473.125
$ws.Range("J80").Value = 1656.7142
$ws.Range("L80").Value = 1656.7142
$ws.Range("N80").Value = -3652.7142
$ws.Range("H83").Value = 1473.125
$ws.Range("J83").Value = 1656.7142
$ws.Range("L83").Value = 8283.571
$ws.Range("N83").Value = -18267.571
$ws.Range("H94").Value = 1122.5
$ws.Range("I94").Value = 1333.3334
$ws.Range("J94").Value = 490
$ws.Range("K94").Value = 1333.3334
$ws.Range("L94").Value = 490
$ws.Range("M94").Value = -882.3334
$ws.Range("N94").Value = -1392
$ws.Range("H133").Value = 49838.75
$ws.Range("J133").Value = 49838.75
$ws.Range("L133").Value = 49838.75
$ws.Range("N133").Value = -59958.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -1000
$ws.Range("H58").Value = 5786.1274
$ws.Range("J58").Value = 6232.0625
$ws.Range("L58").Value = 6232.0625
$ws.Range("N58").Value = -6638.0625
$ws.Range("H86").Value = 8216.75
$ws.Range("J86").Value = 8037.6
$ws.Range("L86").Value = 8037.6
$ws.Range("N86").Value = -10283.6
$ws.Range("H89").Value = 8216.75
$ws.Range("J89").Value = 8037.6
$ws.Range("L89").Value = 40188
$ws.Range("N89").Value = -51420
$ws.Range("H132").Value = 1262.4
$ws.Range("I132").Value = 1271.1666
$ws.Range("J132").Value = 1209.8
$ws.Range("K132").Value = 3813.4998
$ws.Range("L132").Value = 3629.4
$ws.Range("M132").Value = -1283.4998
$ws.Range("N132").Value = -8689.4
$ws.Range("H136").Value = 5786.1274
$ws.Range("J136").Value = 6232.0625
$ws.Range("L136").Value = 18696.1875
$ws.Range("N136").Value = -23796.1875
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 505.2
$ws.Range("I23").Value = 104.666664
$ws.Range("J23").Value = 772.2222
$ws.Range("K23").Value = 313.999992
$ws.Range("L23").Value = 2316.6666
$ws.Range("M23").Value = -78.99999200000002
$ws.Range("N23").Value = -2786.6666
$ws.Range("H33").Value = 3803948.8
$ws.Range("I33").Value = 6173122.5
$ws.Range("K33").Value = 37038735
$ws.Range("M33").Value = -37038452
$ws.Range("H58").Value = 26000
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 6000
$ws.Range("M58").Value = -5872
$ws.Range("H131").Value = 2854.2031
$ws.Range("I131").Value = 1594.25
$ws.Range("J131").Value = 3034.1965
$ws.Range("K131").Value = 4782.75
$ws.Range("L131").Value = 9102.5895
$ws.Range("M131").Value = 257.25
$ws.Range("N131").Value = -19182.5895

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1252771.8
$ws.Range("I80").Value = 836553.8
$ws.Range("J80").Value = 2501425.2
$ws.Range("K80").Value = 836553.8
$ws.Range("L80").Value = 2501425.2
$ws.Range("M80").Value = -835555.8
$ws.Range("N80").Value = -2503421.2
$ws.Range("H83").Value = 1252771.8
$ws.Range("I83").Value = 836553.8
$ws.Range("J83").Value = 2501425.2
$ws.Range("K83").Value = 4182769
$ws.Range("L83").Value = 12507126
$ws.Range("M83").Value = -4177777
$ws.Range("N83").Value = -12517110
$ws.Range("H113").Value = 9433.471
$ws.Range("I113").Value = 2825
$ws.Range("J113").Value = 16868
$ws.Range("K113").Value = 2825
$ws.Range("L113").Value = 16868
$ws.Range("M113").Value = -655
$ws.Range("N113").Value = -21208
$ws.Range("H132").Value = 78829.57000000001
$ws.Range("I132").Value = 9373.429
$ws.Range("K132").Value = 28120.287
$ws.Range("M132").Value = -25590.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1908.3334
$ws.Range("J22").Value = 4500
$ws.Range("L22").Value = 4500
$ws.Range("N22").Value = -5090
$ws.Range("H27").Value = 1908.3334
$ws.Range("J27").Value = 4500
$ws.Range("L27").Value = 4500
$ws.Range("N27").Value = -4714
$ws.Range("H68").Value = 2771.4285
$ws.Range("I68").Value = 2574.5
$ws.Range("J68").Value = 2850.2
$ws.Range("K68").Value = 2574.5
$ws.Range("L68").Value = 2850.2
$ws.Range("M68").Value = -1825.5
$ws.Range("N68").Value = -4348.2
$ws.Range("H71").Value = 2771.4285
$ws.Range("I71").Value = 2574.5
$ws.Range("J71").Value = 2850.2
$ws.Range("K71").Value = 12872.5
$ws.Range("L71").Value = 14251
$ws.Range("M71").Value = -9128.5
$ws.Range("N71").Value = -21739
$ws.Range("H136").Value = 916952
$ws.Range("I136").Value = 1436424.2
$ws.Range("J136").Value = 7875.5
$ws.Range("K136").Value = 4309272.6
$ws.Range("L136").Value = 23626.5
$ws.Range("M136").Value = -4306722.6
$ws.Range("N136").Value = -28726.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 1044
$ws.Range("I42").Value = 1044
$ws.Range("K42").Value = 1044
$ws.Range("M42").Value = -666
$ws.Range("H62").Value = 7111.3335
$ws.Range("I62").Value = 7001
$ws.Range("J62").Value = 7142.857
$ws.Range("K62").Value = 7001
$ws.Range("L62").Value = 7142.857
$ws.Range("M62").Value = -6377
$ws.Range("N62").Value = -8390.857
$ws.Range("H65").Value = 7111.3335
$ws.Range("I65").Value = 7001
$ws.Range("J65").Value = 7142.857
$ws.Range("K65").Value = 35005
$ws.Range("L65").Value = 35714.285
$ws.Range("M65").Value = -31885
$ws.Range("N65").Value = -41954.285
$ws.Range("H96").Value = 145054.72
$ws.Range("I96").Value = 145054.72
$ws.Range("K96").Value = 145054.72
$ws.Range("M96").Value = -143681.72
$ws.Range("H132").Value = 30614.334
$ws.Range("I132").Value = 2148.9033
$ws.Range("K132").Value = 6446.7099
$ws.Range("M132").Value = -3916.7099
